# Adding a new slide.
#
# The original deck has 4 slides:
#   1. Introduction to CSS
#   2. What is CSS
#   3. How to apply CSS to HTML
#   4. Benefits of CSS
#
# A new "Title and Content" slide titled "Adding a new slide" is inserted
# as the new second slide (pushing the rest down).

$p = $ppt.ActivePresentation

# Layout 2 on the slide master is "Title and Content" - the same layout
# used by the existing "What is CSS" slide.
$newSlide = $p.Slides.Add(2, 2)

# Title placeholder (shape 1): type the title text the way PowerPoint
# would when a user types it interactively - "Adding a " followed by
# "new slide" as a second run (e.g. after autocorrect/spell-check splits
# the run on the word boundary).
$titleRange = $newSlide.Shapes.Item(1).TextFrame.TextRange
$titleRange.Text = "Adding a "
$titleRange.InsertAfter("new slide") | Out-Null

# Content placeholder (shape 2) is left empty, as in the authored slide.
